# Scheduled market-data refresh: updates the price/profit columns
# (currentAveragePrice, currentAveragePriceNQ/HQ, LevePriceNQ/HQ,
# LeveProfitNQ/HQ -> columns H-N) for the affected Leve rows on each
# job sheet. Cells that the refresh leaves blank are cleared instead
# of written with 0/empty so the sparse layout round-trips correctly.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 243.61539
$ws.Range("I55").Value = 235.5
$ws.Range("J55").Value = 270.66666
$ws.Range("K55").Value = 235.5
$ws.Range("L55").Value = 270.66666
$ws.Range("M55").Value = -21.5
$ws.Range("N55").Value = -698.66666

$ws.Range("H62").Value = 68464.31
$ws.Range("I62").Value = 115539.445
$ws.Range("J62").Value = 7939.143
$ws.Range("K62").Value = 115539.445
$ws.Range("L62").Value = 7939.143
$ws.Range("M62").Value = -114915.445
$ws.Range("N62").Value = -9187.143

$ws.Range("H65").Value = 68464.31
$ws.Range("I65").Value = 115539.445
$ws.Range("J65").Value = 7939.143
$ws.Range("K65").Value = 577697.2250000001
$ws.Range("L65").Value = 39695.715
$ws.Range("M65").Value = -574577.2250000001
$ws.Range("N65").Value = -45935.715

$ws.Range("H76").Value = 8666.333000000001
$ws.Range("I76").Value = 0
$ws.Range("K76").Value = 0
$ws.Range("M76").ClearContents()

$ws.Range("H79").Value = 8666.333000000001
$ws.Range("I79").Value = 0
$ws.Range("K79").Value = 0
$ws.Range("M79").ClearContents()

$ws.Range("H86").Value = 3485
$ws.Range("I86").Value = 3483
$ws.Range("J86").Value = 3485.5
$ws.Range("K86").Value = 3483
$ws.Range("L86").Value = 3485.5
$ws.Range("M86").Value = -2360
$ws.Range("N86").Value = -5731.5

$ws.Range("H89").Value = 3485
$ws.Range("I89").Value = 3483
$ws.Range("J89").Value = 3485.5
$ws.Range("K89").Value = 17415
$ws.Range("L89").Value = 17427.5
$ws.Range("M89").Value = -11799
$ws.Range("N89").Value = -28659.5

$ws.Range("H106").Value = 2524.5557
$ws.Range("I106").Value = 2102.1428
$ws.Range("K106").Value = 2102.1428
$ws.Range("M106").Value = -1471.1428

$ws.Range("H109").Value = 49500
$ws.Range("J109").Value = 49500
$ws.Range("L109").Value = 49500
$ws.Range("N109").Value = -52274

$ws.Range("H123").Value = 29998
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 29998
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 29998
$ws.Range("M123").ClearContents()
$ws.Range("N123").Value = -39798

$ws.Range("H131").Value = 8211.305
$ws.Range("I131").Value = 1145.1578
$ws.Range("J131").Value = 41775.5
$ws.Range("K131").Value = 3435.4734
$ws.Range("L131").Value = 125326.5
$ws.Range("M131").Value = 1604.5266
$ws.Range("N131").Value = -135406.5

$ws.Range("H132").Value = 2051.6296
$ws.Range("I132").Value = 1726.6364
$ws.Range("K132").Value = 5179.9092
$ws.Range("M132").Value = -2649.9092

$ws.Range("H141").Value = 807.13336
$ws.Range("I141").Value = 807.13336
$ws.Range("K141").Value = 2421.40008
$ws.Range("M141").Value = 2758.59992

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9860.062
$ws.Range("I32").Value = 5745.718
$ws.Range("J32").Value = 25906
$ws.Range("K32").Value = 5745.718
$ws.Range("L32").Value = 25906
$ws.Range("M32").Value = -5458.718
$ws.Range("N32").Value = -26480

$ws.Range("H37").Value = 16000

$ws.Range("H46").Value = 6560.4
$ws.Range("J46").Value = 6289.4443
$ws.Range("L46").Value = 6289.4443
$ws.Range("N46").Value = -6927.4443

$ws.Range("H61").Value = 4368.7446
$ws.Range("I61").Value = 3222.652
$ws.Range("J61").Value = 5467.0835
$ws.Range("K61").Value = 3222.652
$ws.Range("L61").Value = 5467.0835
$ws.Range("M61").Value = -3010.652
$ws.Range("N61").Value = -5891.0835

$ws.Range("H125").Value = 24392
$ws.Range("J125").Value = 24392
$ws.Range("L125").Value = 24392
$ws.Range("N125").Value = -34232

$ws.Range("H132").Value = 5300.5
$ws.Range("I132").Value = 4735.6
$ws.Range("J132").Value = 6995.2
$ws.Range("K132").Value = 14206.8
$ws.Range("L132").Value = 20985.6
$ws.Range("M132").Value = -11676.8
$ws.Range("N132").Value = -26045.6

$ws.Range("H136").Value = 4368.7446
$ws.Range("I136").Value = 3222.652
$ws.Range("J136").Value = 5467.0835
$ws.Range("K136").Value = 9667.956
$ws.Range("L136").Value = 16401.2505
$ws.Range("M136").Value = -7117.956
$ws.Range("N136").Value = -21501.2505

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1126.0333
$ws.Range("I94").Value = 873.2308
$ws.Range("K94").Value = 873.2308
$ws.Range("M94").Value = -422.2308

$ws.Range("H132").Value = 76999.336
$ws.Range("J132").Value = 76999.336
$ws.Range("L132").Value = 76999.336
$ws.Range("N132").Value = -87119.336

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 206.16667
$ws.Range("I22").Value = 206.16667
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 206.16667
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = 143.83333
$ws.Range("N22").ClearContents()

$ws.Range("H31").Value = 35892.734
$ws.Range("I31").Value = 55074.42
$ws.Range("J31").Value = 2760.7273
$ws.Range("K31").Value = 55074.42
$ws.Range("L31").Value = 2760.7273
$ws.Range("M31").Value = -54779.42
$ws.Range("N31").Value = -3350.7273

$ws.Range("H34").Value = 35892.734
$ws.Range("I34").Value = 55074.42
$ws.Range("J34").Value = 2760.7273
$ws.Range("K34").Value = 55074.42
$ws.Range("L34").Value = 2760.7273
$ws.Range("M34").Value = -54872.42
$ws.Range("N34").Value = -3164.7273

$ws.Range("H58").Value = 2335.1614
$ws.Range("I58").Value = 1792.9412
$ws.Range("J58").Value = 2993.5715
$ws.Range("K58").Value = 1792.9412
$ws.Range("L58").Value = 2993.5715
$ws.Range("M58").Value = -1589.9412
$ws.Range("N58").Value = -3399.5715

$ws.Range("H60").Value = 16334
$ws.Range("J60").Value = 14501.5
$ws.Range("L60").Value = 14501.5
$ws.Range("N60").Value = -15523.5

$ws.Range("H86").Value = 1497
$ws.Range("I86").Value = 1497
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 1497
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -374
$ws.Range("N86").ClearContents()

$ws.Range("H89").Value = 1497
$ws.Range("I89").Value = 1497
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 7485
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -1869
$ws.Range("N89").ClearContents()

$ws.Range("H132").Value = 3764.6316
$ws.Range("I132").Value = 3577.4546
$ws.Range("K132").Value = 10732.3638
$ws.Range("M132").Value = -8202.363799999999

$ws.Range("H134").Value = 20724.5
$ws.Range("I134").Value = 7813.6665
$ws.Range("K134").Value = 23440.9995
$ws.Range("M134").Value = -20905.9995

$ws.Range("H136").Value = 2335.1614
$ws.Range("I136").Value = 1792.9412
$ws.Range("J136").Value = 2993.5715
$ws.Range("K136").Value = 5378.8236
$ws.Range("L136").Value = 8980.7145
$ws.Range("M136").Value = -2828.8236
$ws.Range("N136").Value = -14080.7145

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 910.7143
$ws.Range("J113").Value = 945
$ws.Range("L113").Value = 2835
$ws.Range("N113").Value = -7175

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 19042.072
$ws.Range("I70").Value = 15924.5
$ws.Range("J70").Value = 20289.1
$ws.Range("K70").Value = 15924.5
$ws.Range("L70").Value = 20289.1
$ws.Range("M70").Value = -15654.5
$ws.Range("N70").Value = -20829.1

$ws.Range("H73").Value = 19042.072
$ws.Range("I73").Value = 15924.5
$ws.Range("J73").Value = 20289.1
$ws.Range("K73").Value = 15924.5
$ws.Range("L73").Value = 20289.1
$ws.Range("M73").Value = -14988.5
$ws.Range("N73").Value = -22161.1

$ws.Range("H80").Value = 3124.5
$ws.Range("I80").Value = 2999.6667
$ws.Range("K80").Value = 2999.6667
$ws.Range("M80").Value = -2001.6667

$ws.Range("H83").Value = 3124.5
$ws.Range("I83").Value = 2999.6667
$ws.Range("K83").Value = 14998.3335
$ws.Range("M83").Value = -10006.3335

$ws.Range("H132").Value = 4248.923
$ws.Range("I132").Value = 2902.875
$ws.Range("J132").Value = 6402.6
$ws.Range("K132").Value = 8708.625
$ws.Range("L132").Value = 19207.8
$ws.Range("M132").Value = -6178.625
$ws.Range("N132").Value = -24267.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1576.2424
$ws.Range("I16").Value = 1714.7858
$ws.Range("J16").Value = 800.4
$ws.Range("K16").Value = 1714.7858
$ws.Range("L16").Value = 800.4
$ws.Range("M16").Value = -1544.7858
$ws.Range("N16").Value = -1140.4

$ws.Range("H20").Value = 4499.6665
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 4499.6665
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 4499.6665
$ws.Range("M20").ClearContents()
$ws.Range("N20").Value = -4951.6665

$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("M25").ClearContents()
$ws.Range("N25").ClearContents()

$ws.Range("H68").Value = 2963.25
$ws.Range("I68").Value = 3099.4
$ws.Range("J68").Value = 2282.5
$ws.Range("K68").Value = 3099.4
$ws.Range("L68").Value = 2282.5
$ws.Range("M68").Value = -2350.4
$ws.Range("N68").Value = -3780.5

$ws.Range("H71").Value = 2963.25
$ws.Range("I71").Value = 3099.4
$ws.Range("J71").Value = 2282.5
$ws.Range("K71").Value = 15497
$ws.Range("L71").Value = 11412.5
$ws.Range("M71").Value = -11753
$ws.Range("N71").Value = -18900.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 614.7917
$ws.Range("I100").Value = 566.7895
$ws.Range("J100").Value = 797.2
$ws.Range("K100").Value = 1133.579
$ws.Range("L100").Value = 1594.4
$ws.Range("M100").Value = -592.579
$ws.Range("N100").Value = -2676.4

$ws.Range("H122").Value = 2166.7917
$ws.Range("I122").Value = 1732.3125
$ws.Range("K122").Value = 5196.9375
$ws.Range("M122").Value = -2746.9375

